# Daily attendance processing - 2026-02-02 10:20:17 UTC
# Resets all recorded attendance sessions back to "Not Recorded":
#  - clears the "Recorded By" column (G)
#  - zeroes out the recorded-student numerator in "Students" (H), keeping the class size
#  - flips "Status" (I) from Recorded to Not Recorded
#  - updates the roll-up "Class Statistics" (Recorded/Missing Sessions, Coverage %, Avg Attendance %)
#  - updates the "Group Statistics" per-group Recorded/Missing counts and percentages
#  - shrinks/grows a couple of now-unused columns ("Recorded By" / "Status")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 157

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $gCell = $ws.Cells.Item($r, 7)
    $hCell = $ws.Cells.Item($r, 8)
    $iCell = $ws.Cells.Item($r, 9)

    $hVal = $hCell.Value2
    if ($hVal -ne $null -and $hVal -ne "") {
        $parts = [string]$hVal -split '/'
        if ($parts.Length -eq 2) {
            $hCell.Value = "0/" + $parts[1]
        }
    }

    $gCell.Value = ""
    $iCell.Value = "Not Recorded"
}

# --- Class Statistics block ---
# Recorded Sessions (L6) and Missing Sessions (L7) swap: every session is now missing.
$recordedSessions = $ws.Range("L6").Value2
$missingSessions = $ws.Range("L7").Value2
$ws.Range("L6").Value = $missingSessions
$ws.Range("L7").Value = $recordedSessions

# Coverage % and Average Attendance % both drop to 0.0%
$ws.Range("L9").Value = "0.0%"
$ws.Range("L10").Value = "0.0%"

# --- Group Statistics block (rows 15-20): Recorded (O) / Missing (P) swap per group ---
for ($r = 15; $r -le 20; $r++) {
    $oCell = $ws.Cells.Item($r, 15)
    $pCell = $ws.Cells.Item($r, 16)
    $rCell = $ws.Cells.Item($r, 18)
    $sCell = $ws.Cells.Item($r, 19)

    $oVal = $oCell.Value2
    $pVal = $pCell.Value2
    $oCell.Value = $pVal
    $pCell.Value = $oVal

    $rCell.Value = "0.0%"
    $sCell.Value = "0.0%"
}

# --- Column widths: "Recorded By" (G) shrinks now it's empty, "Status" (I) grows a touch ---
$ws.Columns.Item(7).ColumnWidth = 12.140625
$ws.Columns.Item(9).ColumnWidth = 13.140625
